# Refresh the crypto "Price" (column D) and "Volume(1h)" (column E) figures
# in the worksheet, as produced by the scheduled GitHub Actions scraper run.
#
# All of these cells were originally written out as plain text (inline
# strings such as "331.36" or "0.18%"), not numbers/percentages, so we force
# a text ("@") number format on each target cell before assigning the new
# value. That keeps Excel from re-interpreting strings like "45.40" or
# "0.26%" as numeric/percentage values (which would silently drop trailing
# zeros or rescale the percent values).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ref, $value) {
    $range = $ws.Range($ref)
    $range.NumberFormat = "@"
    $range.Value = $value
}

# Row 2 - BNB
Set-TextValue "D2" "331.43"
Set-TextValue "E2" "0.26%"

# Row 3 - OKB
Set-TextValue "D3" "45.43"
Set-TextValue "E3" "2.63%"

# Row 4 - HuobiToken
Set-TextValue "D4" "5.557"
Set-TextValue "E4" "1.09%"

# Row 5 - Cronos
Set-TextValue "D5" "0.08355"
Set-TextValue "E5" "3.96%"

# Row 6 - FTXToken
Set-TextValue "D6" "2.095"
Set-TextValue "E6" "-0.14%"

# Row 7 - MXToken
Set-TextValue "D7" "0.9859"
Set-TextValue "E7" "3.38%"

# Row 8 - BTSEToken
Set-TextValue "D8" "2.546"
Set-TextValue "E8" "-3.69%"

# Row 9 - LiechtensteinCryptoassetsExchange
Set-TextValue "D9" "0.1202"
Set-TextValue "E9" "4.86%"

# Row 10 - WazirX
Set-TextValue "D10" "0.1920"
Set-TextValue "E10" "1.05%"

# Row 11 - MCDex (price unchanged, only volume)
Set-TextValue "E11" "0.58%"

# Row 12 - MandalaExchangeToken
Set-TextValue "D12" "0.09858"
Set-TextValue "E12" "-1.71%"

# Row 13 - BitrueCoin
Set-TextValue "D13" "0.04664"
Set-TextValue "E13" "-3.00%"

# Row 14 - BitMartToken (only price)
Set-TextValue "D14" "0.1060"

# Row 15 - BitForexToken
Set-TextValue "D15" "0.001293"
Set-TextValue "E15" "2.13%"

# Row 16 - TigerCash
Set-TextValue "D16" "0.005946"
Set-TextValue "E16" "2.18%"

# Row 17 - LEO
Set-TextValue "D17" "3.393"
Set-TextValue "E17" "0.77%"

# Row 18 - GateToken
Set-TextValue "D18" "4.449"
Set-TextValue "E18" "0.73%"

# Row 19 - BitpandaEcosystemToken
Set-TextValue "D19" "0.3339"
Set-TextValue "E19" "-3.45%"

# Row 20 - (row 20)
Set-TextValue "D20" "0.1372"
Set-TextValue "E20" "-1.38%"

# Row 21
Set-TextValue "D21" "0.2565"
Set-TextValue "E21" "-0.67%"

# Row 22
Set-TextValue "D22" "0.04153"
Set-TextValue "E22" "1.79%"

# Row 23
Set-TextValue "D23" "0.001294"
Set-TextValue "E23" "1.70%"

# Row 24
Set-TextValue "D24" "0.004565"
Set-TextValue "E24" "4.52%"

# Row 25 (only volume)
Set-TextValue "E25" "8.51%"

# Row 26
Set-TextValue "D26" "0.0003745"
Set-TextValue "E26" "0.03%"

# Row 38 - One
Set-TextValue "D38" "0.02694"
Set-TextValue "E38" "3.83%"

# Row 39 - IDEX
Set-TextValue "D39" "0.05745"
Set-TextValue "E39" "-1.63%"

# Row 40 - KickToken
Set-TextValue "D40" "0.007889"
Set-TextValue "E40" "4.35%"

# Row 41 - BKEXToken
Set-TextValue "D41" "0.1434"
Set-TextValue "E41" "2.20%"

# Row 42 - Dexo
Set-TextValue "D42" "0.007555"
Set-TextValue "E42" "5.35%"

# Row 43 - CEJI
Set-TextValue "D43" "0.002099"
Set-TextValue "E43" "4.15%"

# Row 44 - LocalTraders
Set-TextValue "D44" "0.008937"
Set-TextValue "E44" "8.28%"

# Row 45 - PooCoin (only price)
Set-TextValue "D45" "0.3406"

# Row 46 - CoinLion
Set-TextValue "D46" "0.00007108"
Set-TextValue "E46" "1.19%"

# Row 47 - Kangarootoken
Set-TextValue "D47" "0.00000000751"
Set-TextValue "E47" "0.15%"

# Row 48 - ACDXExchange (only volume)
Set-TextValue "E48" "0.28%"

# Row 49 - BOLO
Set-TextValue "D49" "0.003452"
Set-TextValue "E49" "-1.34%"

# Row 50 - CoinbaseStockToken
Set-TextValue "D50" "0.003534"
Set-TextValue "E50" "0.11%"

# Row 51 - CryptobidCoin
Set-TextValue "D51" "0.00002103"
Set-TextValue "E51" "0.15%"
